$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the next two days of kilométrage (mileage) readings - week 5 data
$ws.Range("A18").Value = 43752
$ws.Range("B18").Value = 439
$ws.Range("A19").Value = 43753
$ws.Range("B19").Value = 480

# Move the active selection to B20, matching where the user left off
$ws.Range("B20").Select()
